$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last row (previously "fold6" / row 7)
$ws.Rows.Item(7).Delete()

# Row 2
$ws.Range("A2").Value = "split1"
$ws.Range("B2").Value = 0.105465445923025
$ws.Range("C2").Value = 0.00968446487245028
$ws.Range("D2").Value = 0.00968446487245028
$ws.Range("E2").Value = 0.0101702486418485
$ws.Range("F2").Value = 0.00690566616132421
$ws.Range("G2").Value = 0.0137561491366818
$ws.Range("H2").Value = 0.12
$ws.Range("I2").Value = 0.117599928568928
$ws.Range("J2").Value = 0.0284239643932029
$ws.Range("K2").Value = 0.000807921751826065
$ws.Range("L2").Value = 0.049221387125471
$ws.Range("M2").Value = 0.208222709856524
$ws.Range("N2").Value = 0.10017008984594
$ws.Range("O2").Value = 0.13959595459099
$ws.Range("P2").Value = 0.0394258647450503
$ws.Range("Q2").Value = 0.170638146864122
$ws.Range("R2").Value = 2.89566232224352
$ws.Range("S2").Value = 0.159001322731053

# Row 3
$ws.Range("A3").Value = "split2"
$ws.Range("B3").Value = 0.102094121770123
$ws.Range("C3").Value = 0.00633447365628936
$ws.Range("D3").Value = 0.00633447365628936
$ws.Range("E3").Value = 0.00418014764034602
$ws.Range("F3").Value = 0.00199638775954267
$ws.Range("G3").Value = 0.00118330607935602
$ws.Range("H3").Value = 0.115
$ws.Range("I3").Value = 0.114611029530331
$ws.Range("J3").Value = 0.0329343323120169
$ws.Range("K3").Value = 0.00108467024483836
$ws.Range("L3").Value = 0.0150737110775756
$ws.Range("M3").Value = 0.213173559349694
$ws.Range("N3").Value = 0.0931662870527073
$ws.Range("O3").Value = 0.137895845610196
$ws.Range("P3").Value = 0.0447295585574887
$ws.Range("Q3").Value = -0.0595052852077916
$ws.Range("R3").Value = 2.77515888691911
$ws.Range("S3").Value = 0.198099848272118

# Row 4
$ws.Range("A4").Value = "split3"
$ws.Range("B4").Value = 0.0958719165605683
$ws.Range("C4").Value = 0.0303518076406184
$ws.Range("D4").Value = 0.0303518076406184
$ws.Range("E4").Value = 0.00119476880739175
$ws.Range("F4").Value = 0.00337695420531697
$ws.Range("G4").Value = 0.00378431244310597
$ws.Range("H4").Value = 0.105
$ws.Range("I4").Value = 0.1076326177571
$ws.Range("J4").Value = 0.037124491469608
$ws.Range("K4").Value = 0.00137822786687699
$ws.Range("L4").Value = 0.00176211874553988
$ws.Range("M4").Value = 0.200674893666289
$ws.Range("N4").Value = 0.0795352931145078
$ws.Range("O4").Value = 0.131079036074453
$ws.Range("P4").Value = 0.0515437429599455
$ws.Range("Q4").Value = -0.149356350223454
$ws.Range("R4").Value = 2.56276266577332
$ws.Range("S4").Value = 0.198912774920749

# Row 5
$ws.Range("A5").Value = "split4"
$ws.Range("B5").Value = 0.119788469302521
$ws.Range("C5").Value = 0.00579689808418174
$ws.Range("D5").Value = 0.00579689808418174
$ws.Range("E5").Value = 0.0026339623431639
$ws.Range("F5").Value = 0.00231796480647126
$ws.Range("G5").Value = 0.00397062082672478
$ws.Range("H5").Value = 0.14
$ws.Range("I5").Value = 0.139827742172983
$ws.Range("J5").Value = 0.0195210508139952
$ws.Range("K5").Value = 0.000381071424882581
$ws.Range("L5").Value = 0.0761262850892554
$ws.Range("M5").Value = 0.19910820854102
$ws.Range("N5").Value = 0.126264592102627
$ws.Range("O5").Value = 0.154583573444775
$ws.Range("P5").Value = 0.0283189813421478
$ws.Range("Q5").Value = -0.0979871391264475
$ws.Range("R5").Value = 2.78186157556275
$ws.Range("S5").Value = 0.122981923451764

# Row 6
$ws.Range("A6").Value = "split5"
$ws.Range("B6").Value = 0.0988019168182415
$ws.Range("C6").Value = 0.0162065226477038
$ws.Range("D6").Value = 0.0162065226477038
$ws.Range("E6").Value = 0.00407762982466713
$ws.Range("F6").Value = 0.00133574258160425
$ws.Range("G6").Value = 0.00168537375991731
$ws.Range("H6").Value = 0.1075
$ws.Range("I6").Value = 0.112944066007418
$ws.Range("J6").Value = 0.0440771493429622
$ws.Range("K6").Value = 0.00194279509420179
$ws.Range("L6").Value = -0.0202213956936644
$ws.Range("M6").Value = 0.225855244061486
$ws.Range("N6").Value = 0.0788523695502407
$ws.Range("O6").Value = 0.138009682055858
$ws.Range("P6").Value = 0.059157312505617
$ws.Range("Q6").Value = -0.29862844077682
$ws.Range("R6").Value = 2.71264764178681
$ws.Range("S6").Value = 0.24607663975515
